# Weekly refresh of the "Fruta, Comercializadora del Agro de Limarí -
# Arándano (blue)" sheet: the Fecha (D), Volumen (M), Precio mínimo (N),
# Precio máximo (O), Precio promedio ponderado (P) and Precio $/Kg (S)
# figures for each market-report row are reshuffled onto different rows
# (dates are stored as Excel day-serials).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  D = 44454; M = 160; N = 12000; O = 13000; P = 12500; S = 6250 },
    @{ Row = 3;  D = 44482; M = 240; N = 10000; O = 11000; P = 10500; S = 5250 },
    @{ Row = 4;  D = 44895; M = 240; N = 3000;  O = 3500;  P = 3250;  S = 1625 },
    @{ Row = 5;  D = 44517; M = 400; N = 5500;  O = 6000;  P = 5750;  S = 2875 },
    @{ Row = 6;  D = 44875; M = 400; N = 7000;  O = 7500;  P = 7250;  S = 3625 },
    @{ Row = 7;  D = 44497; M = 500; N = 9000;  O = 10000; P = 9500;  S = 4750 },
    @{ Row = 8;  D = 44461; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 9;  D = 44819; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 10; D = 44889; M = 460; N = 3500;  O = 4000;  P = 3750;  S = 1875 },
    @{ Row = 11; D = 44489; M = 160; N = 9500;  O = 10000; P = 9750;  S = 4875 },
    @{ Row = 12; D = 44818; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 13; D = 44882; M = 440; N = 6000;  O = 7000;  P = 6500;  S = 3250 },
    @{ Row = 14; D = 44874; M = 300; N = 7500;  O = 8000;  P = 7750;  S = 3875 },
    @{ Row = 15; D = 44475; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 16; D = 44490; M = 400; N = 9500;  O = 10000; P = 9750;  S = 4875 },
    @{ Row = 17; D = 44881; M = 440; N = 6000;  O = 7000;  P = 6500;  S = 3250 },
    @{ Row = 18; D = 44455; M = 200; N = 12000; O = 13000; P = 12500; S = 6250 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # N: Precio mínimo
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O: Precio máximo
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 19).Value = $r.S   # S: Precio $/Kg
}
